# Update the LR-pairs (Wnt1-Fzd3) export with refreshed TPM-derived values.
#
# The sending/target cluster set changed from {FAPs} x {ECs, FAPs, MuSCs,
# Resolving-Mac} to the full {ECs, FAPs, MuSCs} x {ECs, FAPs, MuSCs}
# cross-product (9 data rows instead of 4), and every numeric column was
# recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09370099999999999
$ws.Range("H2").Value = 0.281103
$ws.Range("I2").Value = 0.6634967391997054
$ws.Range("J2").Value = 0.6634967391997054
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.274713
$ws.Range("N2").Value = 0.824139
$ws.Range("O2").Value = 0.1055967877339779
$ws.Range("P2").Value = 0.1055967877339779
$ws.Range("Q2").Value = 0.025740882813
$ws.Range("R2").Value = 0.231667945317
$ws.Range("S2").Value = 0.07006312433145777
$ws.Range("T2").Value = 0.07006312433145777

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09370099999999999
$ws.Range("H3").Value = 0.281103
$ws.Range("I3").Value = 0.6634967391997054
$ws.Range("J3").Value = 0.6634967391997054
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8886716666666666
$ws.Range("N3").Value = 2.666015
$ws.Range("O3").Value = 0.3415960415058637
$ws.Range("P3").Value = 0.3415960415058638
$ws.Range("Q3").Value = 0.08326942383833333
$ws.Range("R3").Value = 0.7494248145449999
$ws.Range("S3").Value = 0.2266478596626678
$ws.Range("T3").Value = 0.2266478596626678

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09370099999999999
$ws.Range("H4").Value = 0.281103
$ws.Range("I4").Value = 0.6634967391997054
$ws.Range("J4").Value = 0.6634967391997054
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.438143333333333
$ws.Range("N4").Value = 4.31443
$ws.Range("O4").Value = 0.5528071707601584
$ws.Range("P4").Value = 0.5528071707601584
$ws.Range("Q4").Value = 0.1347554684766667
$ws.Range("R4").Value = 1.21279921629
$ws.Range("S4").Value = 0.3667857552055798
$ws.Range("T4").Value = 0.3667857552055798

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.023296
$ws.Range("H5").Value = 0.069888
$ws.Range("I5").Value = 0.1649589656075852
$ws.Range("J5").Value = 0.1649589656075852
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.274713
$ws.Range("N5").Value = 0.824139
$ws.Range("O5").Value = 0.1055967877339779
$ws.Range("P5").Value = 0.1055967877339779
$ws.Range("Q5").Value = 0.006399714048
$ws.Range("R5").Value = 0.057597426432
$ws.Range("S5").Value = 0.01741913687608073
$ws.Range("T5").Value = 0.01741913687608073

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.023296
$ws.Range("H6").Value = 0.069888
$ws.Range("I6").Value = 0.1649589656075852
$ws.Range("J6").Value = 0.1649589656075852
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8886716666666666
$ws.Range("N6").Value = 2.666015
$ws.Range("O6").Value = 0.3415960415058637
$ws.Range("P6").Value = 0.3415960415058638
$ws.Range("Q6").Value = 0.02070249514666667
$ws.Range("R6").Value = 0.18632245632
$ws.Range("S6").Value = 0.05634932966245301
$ws.Range("T6").Value = 0.05634932966245301

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.023296
$ws.Range("H7").Value = 0.069888
$ws.Range("I7").Value = 0.1649589656075852
$ws.Range("J7").Value = 0.1649589656075852
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.438143333333333
$ws.Range("N7").Value = 4.31443
$ws.Range("O7").Value = 0.5528071707601584
$ws.Range("P7").Value = 0.5528071707601584
$ws.Range("Q7").Value = 0.03350298709333333
$ws.Range("R7").Value = 0.30152688384
$ws.Range("S7").Value = 0.09119049906905144
$ws.Range("T7").Value = 0.09119049906905143

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.024226
$ws.Range("H8").Value = 0.072678
$ws.Range("I8").Value = 0.1715442951927094
$ws.Range("J8").Value = 0.1715442951927094
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.274713
$ws.Range("N8").Value = 0.824139
$ws.Range("O8").Value = 0.1055967877339779
$ws.Range("P8").Value = 0.1055967877339779
$ws.Range("Q8").Value = 0.006655197138
$ws.Range("R8").Value = 0.059896774242
$ws.Range("S8").Value = 0.01811452652643938
$ws.Range("T8").Value = 0.01811452652643938

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.024226
$ws.Range("H9").Value = 0.072678
$ws.Range("I9").Value = 0.1715442951927094
$ws.Range("J9").Value = 0.1715442951927094
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8886716666666666
$ws.Range("N9").Value = 2.666015
$ws.Range("O9").Value = 0.3415960415058637
$ws.Range("P9").Value = 0.3415960415058638
$ws.Range("Q9").Value = 0.02152895979666667
$ws.Range("R9").Value = 0.19376063817
$ws.Range("S9").Value = 0.0585988521807429
$ws.Range("T9").Value = 0.05859885218074291

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt1"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.024226
$ws.Range("H10").Value = 0.072678
$ws.Range("I10").Value = 0.1715442951927094
$ws.Range("J10").Value = 0.1715442951927094
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.438143333333333
$ws.Range("N10").Value = 4.31443
$ws.Range("O10").Value = 0.5528071707601584
$ws.Range("P10").Value = 0.5528071707601584
$ws.Range("Q10").Value = 0.03484046039333333
$ws.Range("R10").Value = 0.31356414354
$ws.Range("S10").Value = 0.09483091648552713
$ws.Range("T10").Value = 0.09483091648552713

